$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the comment that was attached to A3, then clear the cell itself so
# the whole row 3 disappears from the sheet (row 5 keeps its own row number).
$ws.Range("A3").Comment.Delete()
$ws.Range("A3").ClearContents()

# Give A5 a thin box border on all four sides (default/automatic color).
$ws.Range("A5").Borders.LineStyle = 1
$ws.Range("A5").Borders.Weight = 2
$ws.Range("A5").Borders.ColorIndex = -4105

# Move the active selection to F4.
$ws.Range("F4").Select()
